# Weekly update: insert a new price record as row 194 in the Papa
# (Terminal Hortofrutícola Agro Chillán) sheet, pushing the existing
# rows 194-256 down to 195-257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194 (shifts rows 194:256 -> 195:257)
$ws.Rows(194).Insert()

# Populate the newly inserted row with the latest week's data
$ws.Range("A194").Value = 7
$ws.Range("B194").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C194").Value = "Ñuble"
$ws.Range("D194").Value = 44559
$ws.Range("E194").Value = 16
$ws.Range("F194").Value = 100114001
$ws.Range("G194").Value = "Papa"
$ws.Range("H194").Value = "Asterix"
$ws.Range("I194").Value = "1a nueva(o)"
$ws.Range("J194").Value = 160
$ws.Range("K194").Value = 9500
$ws.Range("L194").Value = 10000
$ws.Range("M194").Value = 9750
$ws.Range("N194").Value = "`$/saco 25 kilos"
$ws.Range("O194").Value = "Región del Maule"
$ws.Range("P194").Value = 390
$ws.Range("Q194").Value = 25
$ws.Range("R194").Value = "Hortaliza"
